$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

$ws.Range("D2").Value = 0.003099348396062851
$ws.Range("E2").Value = 1.461850028019398
$ws.Range("G2").Value = 0.09455263894051313
$ws.Range("H2").Value = 0.7239843839779496
$ws.Range("I2").Value = 0.2240070682018995
$ws.Range("J2").Value = 0.3059666785411537
$ws.Range("K2").Value = 0.02817187272012234

$ws.Range("D3").Value = 0.1007486316375434
$ws.Range("E3").Value = 0.6844639610499144
$ws.Range("G3").Value = 0.0271370792761445
$ws.Range("H3").Value = 0.31308091012761
$ws.Range("I3").Value = 0.03509241668507457
$ws.Range("J3").Value = 0.2774711814709008
$ws.Range("K3").Value = 0.007966393604874611

$ws.Range("C4").Value = 865
$ws.Range("D4").Value = 0.1145646297372878
$ws.Range("E4").Value = 0.8280656086280942
$ws.Range("F4").Value = 865
$ws.Range("G4").Value = 0.03680296847596765
$ws.Range("H4").Value = 0.4023053646087646
$ws.Range("I4").Value = 0.03863512258976698
$ws.Range("J4").Value = 0.3080506366677582
$ws.Range("K4").Value = 0.01047915313392878

$ws.Range("D5").Value = 0.004071842413395643
$ws.Range("E5").Value = 1.445242314599454
$ws.Range("G5").Value = 0.09286676626652479
$ws.Range("H5").Value = 0.7123333448544145
$ws.Range("I5").Value = 0.2266222876496613
$ws.Range("J5").Value = 0.3017501258291304
$ws.Range("K5").Value = 0.02738049998879433

$ws.Range("D6").Value = 0.1413284973241389
$ws.Range("E6").Value = 1.144778670743108
$ws.Range("G6").Value = 0.04122060397639871
$ws.Range("H6").Value = 0.4348260699771345
$ws.Range("I6").Value = 0.3393531036563218
$ws.Range("J6").Value = 0.2823181990534067
$ws.Range("K6").Value = 0.01322491047903895

$ws.Range("E7").Value = 23.88862068532035

$ws.Range("D8").Value = 0.003069788217544556
$ws.Range("E8").Value = 1.546455363743007
$ws.Range("G8").Value = 0.09828882524743676
$ws.Range("H8").Value = 0.7944544311612844
$ws.Range("I8").Value = 0.2120680687949061
$ws.Range("J8").Value = 0.3219255269505084
$ws.Range("K8").Value = 0.0294781387783587

$ws.Range("D9").Value = 0.1306328712962568
$ws.Range("E9").Value = 0.9219073071144521
$ws.Range("G9").Value = 0.04106759186834097
$ws.Range("H9").Value = 0.4729144708253443
$ws.Range("I9").Value = 0.04559832625091076
$ws.Range("J9").Value = 0.3105810107663274
$ws.Range("K9").Value = 0.01279717171564698

$ws.Range("C10").Value = 3578
$ws.Range("D10").Value = 0.1551502710208297
$ws.Range("E10").Value = 2.179413402918726
$ws.Range("F10").Value = 3578
$ws.Range("G10").Value = 0.1384429661557078
$ws.Range("H10").Value = 1.287329831160605
$ws.Range("I10").Value = 0.0587250916287303
$ws.Range("J10").Value = 0.5194771252572536
$ws.Range("K10").Value = 0.04343154653906822

$ws.Range("D11").Value = 0.004328818060457706
$ws.Range("E11").Value = 1.487140614073724
$ws.Range("G11").Value = 0.09522184357047081
$ws.Range("H11").Value = 0.7370280637405813
$ws.Range("I11").Value = 0.2154956161975861
$ws.Range("J11").Value = 0.3227867158129811
$ws.Range("K11").Value = 0.02909201802685857

$ws.Range("D12").Value = 0.1795502840541303
$ws.Range("E12").Value = 1.780028450768441
$ws.Range("G12").Value = 0.06949109956622124
$ws.Range("H12").Value = 0.6853195149451494
$ws.Range("I12").Value = 0.6357602667994797
$ws.Range("J12").Value = 0.3061654586344957
$ws.Range("K12").Value = 0.02337892353534698

$ws.Range("E13").Value = 24.3470804435201

$ws.Range("D14").Value = 0.003084568306803703
$ws.Range("E14").Value = 1.504152695881203
$ws.Range("G14").Value = 0.09642073209397495
$ws.Range("H14").Value = 0.759219407569617
$ws.Range("I14").Value = 0.2180375684984028
$ws.Range("J14").Value = 0.313946102745831
$ws.Range("K14").Value = 0.02882500574924052

$ws.Range("D15").Value = 0.1156907514669001
$ws.Range("E15").Value = 0.8031856340821832
$ws.Range("G15").Value = 0.03410233557224274
$ws.Range("H15").Value = 0.3929976904764771
$ws.Range("I15").Value = 0.04034537146799266
$ws.Range("J15").Value = 0.2940260961186141
$ws.Range("K15").Value = 0.0103817826602608

$ws.Range("C16").Value = 2221.5
$ws.Range("D16").Value = 0.1348574503790587
$ws.Range("E16").Value = 1.50373950577341
$ws.Range("F16").Value = 2221.5
$ws.Range("G16").Value = 0.08762296731583774
$ws.Range("H16").Value = 0.8448175978846848
$ws.Range("I16").Value = 0.04868010710924864
$ws.Range("J16").Value = 0.4137638809625059
$ws.Range("K16").Value = 0.0269553498364985

$ws.Range("D17").Value = 0.004200330236926675
$ws.Range("E17").Value = 1.466191464336589
$ws.Range("G17").Value = 0.0940443049184978
$ws.Range("H17").Value = 0.7246807042974979
$ws.Range("I17").Value = 0.2210589519236237
$ws.Range("J17").Value = 0.3122684208210558
$ws.Range("K17").Value = 0.02823625900782645

$ws.Range("D18").Value = 0.1604393906891346
$ws.Range("E18").Value = 1.462403560755774
$ws.Range("G18").Value = 0.05535585177130997
$ws.Range("H18").Value = 0.5600727924611419
$ws.Range("I18").Value = 0.4875566852279007
$ws.Range("J18").Value = 0.2942418288439512
$ws.Range("K18").Value = 0.01830191700719297

$ws.Range("E19").Value = 24.11785056442022
